$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.485259333333333
$ws.Range("N2").Value = 4.455778
$ws.Range("O2").Value = 0.3057455162066235
$ws.Range("P2").Value = 0.3057455162066235
$ws.Range("Q2").Value = 2.554719821213556
$ws.Range("R2").Value = 22.99247839092201
$ws.Range("S2").Value = 0.1971784018447014
$ws.Range("T2").Value = 0.1971784018447014

$ws.Range("O3").Value = 0.2805555239151429
$ws.Range("P3").Value = 0.2805555239151429
$ws.Range("S3").Value = 0.1809331188912865
$ws.Range("T3").Value = 0.1809331188912865

$ws.Range("O4").Value = 0.4136989598782336
$ws.Range("P4").Value = 0.4136989598782336
$ws.Range("S4").Value = 0.2667986787367258
$ws.Range("T4").Value = 0.2667986787367258

$ws.Range("M5").Value = 1.485259333333333
$ws.Range("N5").Value = 4.455778
$ws.Range("O5").Value = 0.3057455162066235
$ws.Range("P5").Value = 0.3057455162066235
$ws.Range("Q5").Value = 1.406637625609778
$ws.Range("R5").Value = 12.659738630488
$ws.Range("S5").Value = 0.1085671143619221
$ws.Range("T5").Value = 0.1085671143619221

$ws.Range("O6").Value = 0.2805555239151429
$ws.Range("P6").Value = 0.2805555239151429
$ws.Range("S6").Value = 0.0996224050238564
$ws.Range("T6").Value = 0.09962240502385641

$ws.Range("O7").Value = 0.4136989598782336
$ws.Range("P7").Value = 0.4136989598782336
$ws.Range("S7").Value = 0.1469002811415078
$ws.Range("T7").Value = 0.1469002811415078
